# Auto-generated Excel COM-interop script applying the scheduled market-data refresh
# (currentAveragePrice* / LevePrice* / LeveProfit* columns) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 253.4
$ws.Range("I2").Value = 241.75
$ws.Range("K2").Value = 241.75
$ws.Range("M2").Value = -128.75
$ws.Range("H15").Value = 1346.0676
$ws.Range("I15").Value = 1346.0676
$ws.Range("K15").Value = 4038.2028
$ws.Range("M15").Value = -3869.2028
$ws.Range("H55").Value = 166668510
$ws.Range("I55").Value = 500005000
$ws.Range("J55").Value = 275
$ws.Range("K55").Value = 500005000
$ws.Range("L55").Value = 275
$ws.Range("M55").Value = -500004786
$ws.Range("N55").Value = -703
$ws.Range("H70").Value = 907.0769
$ws.Range("I70").Value = 817.3333
$ws.Range("J70").Value = 984
$ws.Range("K70").Value = 2451.9999
$ws.Range("L70").Value = 2952
$ws.Range("M70").Value = -2181.9999
$ws.Range("N70").Value = -3492
$ws.Range("H73").Value = 907.0769
$ws.Range("I73").Value = 817.3333
$ws.Range("J73").Value = 984
$ws.Range("K73").Value = 2451.9999
$ws.Range("L73").Value = 2952
$ws.Range("M73").Value = -1515.9999
$ws.Range("N73").Value = -4824
$ws.Range("H87").Value = 24103.393
$ws.Range("J87").Value = 24103.393
$ws.Range("L87").Value = 24103.393
$ws.Range("N87").Value = -26599.393
$ws.Range("H90").Value = 24103.393
$ws.Range("J90").Value = 24103.393
$ws.Range("L90").Value = 72310.179
$ws.Range("N90").Value = -84790.179
$ws.Range("H113").Value = 2565.8333
$ws.Range("J113").Value = 2598.75
$ws.Range("L113").Value = 2598.75
$ws.Range("N113").Value = -9106.75
$ws.Range("H116").Value = 2554.5833
$ws.Range("I116").Value = 2422.7778
$ws.Range("K116").Value = 2422.7778
$ws.Range("M116").Value = 1019.2222
$ws.Range("H137").Value = 8584.166999999999
$ws.Range("I137").Value = 9800.4
$ws.Range("K137").Value = 29401.2
$ws.Range("M137").Value = -26851.2
$ws.Range("H138").Value = 193171.45
$ws.Range("I138").Value = 2049.0435
$ws.Range("J138").Value = 322460.16
$ws.Range("K138").Value = 6147.130500000001
$ws.Range("L138").Value = 967380.48
$ws.Range("M138").Value = -1007.130500000001
$ws.Range("N138").Value = -977660.48
$ws.Range("H141").Value = 3447.318
$ws.Range("I141").Value = 1648.2759
$ws.Range("K141").Value = 4944.8277
$ws.Range("M141").Value = 235.1723000000002

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H45").Value = 3333.2666
$ws.Range("I45").Value = 2642.7144
$ws.Range("K45").Value = 2642.7144
$ws.Range("M45").Value = -2265.7144
$ws.Range("H74").Value = 1389
$ws.Range("I74").Value = 1066.8
$ws.Range("J74").Value = 3000
$ws.Range("K74").Value = 1066.8
$ws.Range("L74").Value = 3000
$ws.Range("M74").Value = -192.8
$ws.Range("N74").Value = -4748
$ws.Range("H77").Value = 1389
$ws.Range("I77").Value = 1066.8
$ws.Range("J77").Value = 3000
$ws.Range("K77").Value = 5334
$ws.Range("L77").Value = 15000
$ws.Range("M77").Value = -966
$ws.Range("N77").Value = -23736
$ws.Range("H132").Value = 5264.2
$ws.Range("I132").Value = 4678.0713
$ws.Range("K132").Value = 14034.2139
$ws.Range("M132").Value = -11504.2139

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H20").Value = 1705.35
$ws.Range("I20").Value = 1479.1428
$ws.Range("J20").Value = 2233.1667
$ws.Range("K20").Value = 1479.1428
$ws.Range("L20").Value = 2233.1667
$ws.Range("M20").Value = -1232.1428
$ws.Range("N20").Value = -2727.1667
$ws.Range("H94").Value = 528.1875
$ws.Range("I94").Value = 458.36365
$ws.Range("J94").Value = 681.8
$ws.Range("K94").Value = 458.36365
$ws.Range("L94").Value = 681.8
$ws.Range("M94").Value = -7.363650000000007
$ws.Range("N94").Value = -1583.8
$ws.Range("H134").Value = 4235.9414
$ws.Range("I134").Value = 4401.375
$ws.Range("J134").Value = 4088.889
$ws.Range("K134").Value = 13204.125
$ws.Range("L134").Value = 12266.667
$ws.Range("M134").Value = -10669.125
$ws.Range("N134").Value = -17336.667

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5356.375
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 5356.375
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 5356.375
$ws.Range("N31").Value = -5946.375
$ws.Range("M31").ClearContents()
$ws.Range("H34").Value = 5356.375
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 5356.375
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 5356.375
$ws.Range("N34").Value = -5760.375
$ws.Range("M34").ClearContents()

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()
$ws.Range("H132").Value = 3744.0344
$ws.Range("J132").Value = 4407.8286
$ws.Range("L132").Value = 39670.4574
$ws.Range("N132").Value = -44730.4574

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1799.9333
$ws.Range("I97").Value = 1778.8889
$ws.Range("J97").Value = 1831.5
$ws.Range("K97").Value = 1778.8889
$ws.Range("L97").Value = 1831.5
$ws.Range("M97").Value = -1282.8889
$ws.Range("N97").Value = -2823.5
$ws.Range("H122").Value = 4507.1562
$ws.Range("I122").Value = 3294.1
$ws.Range("J122").Value = 5058.5454
$ws.Range("K122").Value = 9882.299999999999
$ws.Range("L122").Value = 15175.6362
$ws.Range("M122").Value = -7432.299999999999
$ws.Range("N122").Value = -20075.6362
$ws.Range("H132").Value = 2825.1
$ws.Range("I132").Value = 2208.6667
$ws.Range("J132").Value = 3749.75
$ws.Range("K132").Value = 6626.000100000001
$ws.Range("L132").Value = 11249.25
$ws.Range("M132").Value = -4096.000100000001
$ws.Range("N132").Value = -16309.25
$ws.Range("H133").Value = 55000
$ws.Range("J133").Value = 55000
$ws.Range("L133").Value = 55000
$ws.Range("N133").Value = -65120

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 83336504
$ws.Range("I7").Value = 100003320
$ws.Range("J7").Value = 2402.5
$ws.Range("K7").Value = 100003320
$ws.Range("L7").Value = 2402.5
$ws.Range("M7").Value = -100003208
$ws.Range("N7").Value = -2626.5
$ws.Range("H46").Value = 1224.75
$ws.Range("I46").Value = 1350
$ws.Range("J46").Value = 1099.5
$ws.Range("K46").Value = 1350
$ws.Range("L46").Value = 1099.5
$ws.Range("M46").Value = -1162
$ws.Range("N46").Value = -1475.5
$ws.Range("H126").Value = 83336504
$ws.Range("I126").Value = 100003320
$ws.Range("J126").Value = 2402.5
$ws.Range("K126").Value = 300009960
$ws.Range("L126").Value = 7207.5
$ws.Range("M126").Value = -300007490
$ws.Range("N126").Value = -12147.5
$ws.Range("H132").Value = 3892.325
$ws.Range("I132").Value = 3564.64
$ws.Range("J132").Value = 4438.467
$ws.Range("K132").Value = 10693.92
$ws.Range("L132").Value = 13315.401
$ws.Range("M132").Value = -8163.92
$ws.Range("N132").Value = -18375.401
$ws.Range("H133").Value = 34975.5
$ws.Range("J133").Value = 34975.5
$ws.Range("L133").Value = 34975.5
$ws.Range("N133").Value = -40035.5
$ws.Range("H136").Value = 9261177
$ws.Range("I136").Value = 2396
$ws.Range("K136").Value = 7188
$ws.Range("M136").Value = -4638

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1735.5714
$ws.Range("I126").Value = 1350
$ws.Range("J126").Value = 1889.8
$ws.Range("K126").Value = 4050
$ws.Range("L126").Value = 5669.4
$ws.Range("M126").Value = -1580
$ws.Range("N126").Value = -10609.4
$ws.Range("H132").Value = 15157235
$ws.Range("I132").Value = 13335
$ws.Range("J132").Value = 20836196
$ws.Range("K132").Value = 40005
$ws.Range("L132").Value = 62508588
$ws.Range("M132").Value = -37475
$ws.Range("N132").Value = -62513648
$ws.Range("H136").Value = 3505.6667
$ws.Range("I136").Value = 3246.0527
$ws.Range("J136").Value = 3954.0908
$ws.Range("K136").Value = 9738.158100000001
$ws.Range("L136").Value = 11862.2724
$ws.Range("M136").Value = -9188.158100000001
$ws.Range("N136").Value = -16962.2724
